# Fix unicode errors and improve text file formatting:
# strip stray trailing spaces that were left at the end of several
# lyric lines (runs followed by a line break) on the Verse 1 (slide 1),
# Chorus (slide 2) and Bridge (slide 4) text boxes.
#
# We edit via TextRange.Characters(start, length) so only the run's
# <a:t> content is touched - the surrounding <a:r>/<a:br/> structure is
# left completely intact, matching the target diff.

function Replace-RunText {
    param($TextRange, $OldText, $NewText, $FromIndex)
    $full = $TextRange.Text
    $idx = $full.IndexOf($OldText, $FromIndex)
    if ($idx -lt 0) {
        throw "Text not found: [$OldText] starting at $FromIndex"
    }
    $sub = $TextRange.Characters($idx + 1, $OldText.Length)
    $sub.Text = $NewText
    return $idx + $NewText.Length
}

$p = $ppt.ActivePresentation

# --- Slide 1 (Verse 1) ---
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(1).TextFrame.TextRange
$pos = 0
$pos = Replace-RunText $tr1 "In a dry and barren land I bow down  " "In a dry and barren land I bow down" $pos
$pos = Replace-RunText $tr1 "I need You now " "I need You now" $pos
$pos = Replace-RunText $tr1 "To Your river I will run I bow down  " "To Your river I will run I bow down" $pos
$pos = Replace-RunText $tr1 "I need You now " "I need You now" $pos

# --- Slide 2 (Chorus) ---
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(1).TextFrame.TextRange
$pos = 0
$pos = Replace-RunText $tr2 "O Living Water O God my Savior  " "O Living Water O God my Savior" $pos
$pos = Replace-RunText $tr2 "If I ever needed You I need You now  " "If I ever needed You I need You now" $pos
$pos = Replace-RunText $tr2 "O Living Water O God my Healer  " "O Living Water O God my Healer" $pos
$pos = Replace-RunText $tr2 "If I ever needed You I need You now  " "If I ever needed You I need You now" $pos

# --- Slide 4 (Bridge) ---
$s4 = $p.Slides.Item(4)
$tr4 = $s4.Shapes.Item(1).TextFrame.TextRange
$pos = 0
$pos = Replace-RunText $tr4 "Just like the desert needs the blessing of the rain  " "Just like the desert needs the blessing of the rain" $pos
$pos = Replace-RunText $tr4 "Just like the winter waiting for the sun again  " "Just like the winter waiting for the sun again" $pos
$pos = Replace-RunText $tr4 "I need You now  " "I need You now" $pos
